$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8..103 down to 9..104
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new record
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "Macroferia Regional de Talca"
$ws.Range("C8").Value = "Maule"
$ws.Range("D8").Value = 45190
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 300000000
$ws.Range("G8").Value = "Espárragos"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 1500
$ws.Range("K8").Value = 1200
$ws.Range("L8").Value = 1200
$ws.Range("M8").Value = 1200
$ws.Range("N8").Value = "$/kilo"
$ws.Range("O8").Value = "Provincia de Linares"
$ws.Range("P8").Value = 1200
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = "Hortaliza"
